# "Added last minute updates" — AF PGI 5315 cover paragraph tweaks:
#   * replace the merge-field placeholder text
#   * drop the now-unused trailing space run
#   * tighten the left indent
#   * add a (invisible, space-only) paragraph border on all four sides

$d = $word.ActiveDocument

# --- First paragraph: update the bookmark-style placeholder text ---
$d.Content.Find.Execute("**ID__AFFARS_pgi_5315_topic_2__ID**", $true, $false, $false, $false, $false,
                         $true, 1, $false, "**ID__AFFARS_AF_PGI_5315__ID**", 2) | Out-Null

$p1 = $d.Paragraphs(1)

# --- Remove the now-orphaned trailing space run at the end of the paragraph ---
$pRange = $p1.Range
$trailingSpace = $d.Range($pRange.End - 2, $pRange.End - 1)
if ($trailingSpace.Text -eq " ") {
    $trailingSpace.Delete()
}

# --- Update paragraph indentation (120 -> 225 twips, i.e. 9pt -> 11.25pt) ---
$p1.Format.LeftIndent = 11.25

# --- Add a paragraph border (space-only, no visible line) on all four sides ---
$p1.Format.Borders.DistanceFromTop = 5
$p1.Format.Borders.DistanceFromLeft = 5
$p1.Format.Borders.DistanceFromBottom = 5
$p1.Format.Borders.DistanceFromRight = 5
